# Update scores for evolution algorithm.
# The evolution algorithm was revised, so the Z column scores (rows 3-22)
# on Sheet1 are updated with new test results. Z23 holds
# =AVERAGE(Z3:Z22) and recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Z3").Value  = 3040
$ws.Range("Z4").Value  = 2350
$ws.Range("Z5").Value  = 6150
$ws.Range("Z6").Value  = 3360
$ws.Range("Z7").Value  = 1630
$ws.Range("Z8").Value  = 2710
$ws.Range("Z9").Value  = 5580
$ws.Range("Z10").Value = 2720
$ws.Range("Z11").Value = 4940
$ws.Range("Z12").Value = 960
$ws.Range("Z13").Value = 1120
$ws.Range("Z14").Value = 2860
$ws.Range("Z15").Value = 2710
$ws.Range("Z16").Value = 1650
$ws.Range("Z17").Value = 4250
$ws.Range("Z18").Value = 5020
$ws.Range("Z19").Value = 2570
$ws.Range("Z20").Value = 710
$ws.Range("Z21").Value = 2340
$ws.Range("Z22").Value = 1610

# Move the active selection to Z23, matching the last recorded cursor
# position after the data was revised.
[void]$ws.Range("Z23").Select()
